# xlsx-exporter test workbook: add "arr1"/"arr2" demo columns to the
# "follow" sheet so the expr-checker can be exercised against array-typed
# columns (commit: "use '$$' to get row data in expr checker").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("follow")

# -- header / type / marker rows -------------------------------------------------
$ws.Range("D1").Value = "arr1"
$ws.Range("E1").Value = "arr2"

$ws.Range("D2").Value = "int[]?"
$ws.Range("E2").Value = "int[]?"

# row 3 ('>>') and row 5 ('###') have nothing in the new columns - leave blank

# row 6 ('!!!' -> first sample row) carries the sample array literal; fill it
# in before the validation-expression cell so the shared-string table order
# matches the authoring order.
$ws.Range("D6").Value = "[1,2]"
$ws.Range("E6").Value = "[1,2]"

# row 4 ('!!!' marker row): plain "x" for arr1, and the new expr-checker
# formula for arr2 that cross-references the sibling "arr1" column.
$ws.Range("D4").Value = "x"
$ws.Range("E4").Value = "`$.length == arr1.length"

# -- formatting for the new data column cells (rows 6-13) ------------------------
$dataRange = $ws.Range("D6:E13")
$dataRange.Borders.LineStyle = 1
$dataRange.VerticalAlignment = -4160

# make this sheet the active one, matching the saved workbook's active tab
$ws.Activate()
